# Implementação da edição de endereço
# Implementação do parametros para identificação do relatório no gerador de relatório de contas a receber
#
# Adds 7 new backlog rows (32-38) to the "Itens de desenvolvimento" sheet,
# reusing the existing "Melhoria/Defeito/Desenvolvido/N/A" shared strings
# and matching the look (fill/border) of the neighbouring rows, then moves
# the active selection to reflect where the user ended up after typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Two "Defeito" rows, styled like row 31 (plain green-ish list rows)
# ---------------------------------------------------------------------
$ws.Range("A31:D31").Copy($ws.Range("A32:D32"))
$ws.Cells.Item(32, 1).Value = "Retirada do campo obervação no relatório de boletos gerados "
$ws.Cells.Item(32, 2).Value = "Defeito"
$ws.Cells.Item(32, 3).Value = "Desenvolvido"
$ws.Cells.Item(32, 4).Value = "N/A"

$ws.Range("A31:D31").Copy($ws.Range("A33:D33"))
$ws.Cells.Item(33, 1).Value = "Correção na quantidade limite de informações carregadas nos relatório de contas a receber e boletos gerados"
$ws.Cells.Item(33, 2).Value = "Defeito"
$ws.Cells.Item(33, 3).Value = "Desenvolvido"
$ws.Cells.Item(33, 4).Value = "N/A"
$ws.Rows.Item(33).RowHeight = 25.5

# ---------------------------------------------------------------------
# 2) Five "Melhoria" rows, styled like row 29 (yellow highlighted rows)
# ---------------------------------------------------------------------
$ws.Range("A29:D29").Copy($ws.Range("A34:D34"))
$ws.Cells.Item(34, 1).Value = "Melhoria no leiaute do boleto"
$ws.Cells.Item(34, 2).Value = "Melhoria"
$ws.Cells.Item(34, 3).Value = "Desenvolvido"
$ws.Cells.Item(34, 4).Value = "N/A"

$ws.Range("A29:D29").Copy($ws.Range("A35:D35"))
$ws.Cells.Item(35, 1).Value = "Permitir gerar itens financeiros e posteriormente gerar as contas a receber "
$ws.Cells.Item(35, 2).Value = "Melhoria"
$ws.Cells.Item(35, 3).Value = "Desenvolvido"
$ws.Cells.Item(35, 4).Value = "N/A"

$ws.Range("A29:D29").Copy($ws.Range("A36:D36"))
$ws.Cells.Item(36, 1).Value = "Permitir cadastrar mais de um e-mail para uma pessoa e poder editá-los"
$ws.Cells.Item(36, 2).Value = "Melhoria"
$ws.Cells.Item(36, 3).Value = "Desenvolvido"
$ws.Cells.Item(36, 4).Value = "N/A"

$ws.Range("A29:D29").Copy($ws.Range("A37:D37"))
$ws.Cells.Item(37, 1).Value = "Criar opção que permite cadastrar um contato para cada telefone da pessoa"
$ws.Cells.Item(37, 2).Value = "Melhoria"
$ws.Cells.Item(37, 3).Value = "Desenvolvido"
$ws.Cells.Item(37, 4).Value = "N/A"

$ws.Range("A29:D29").Copy($ws.Range("A38:D38"))
$ws.Cells.Item(38, 1).Value = "Permitir editar os endereços das pessoas"
$ws.Cells.Item(38, 2).Value = "Melhoria"
$ws.Cells.Item(38, 3).Value = "Desenvolvido"
$ws.Cells.Item(38, 4).Value = "N/A"

# ---------------------------------------------------------------------
# 3) Leave the workbook where the author left it: scrolled down with the
#    next empty row selected, ready for further data entry.
# ---------------------------------------------------------------------
$ws.Range("A15").Select()
$ws.Range("A39").Select()
